# Add a new row of credentials data (row 4) to the "Credentials" sheet,
# matching the layout/styling of the existing rows (2 and 3), and update
# the active selection to D9.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- A4: "valid" (same plain text style as A2/A3) ---
$ws.Range("A2").Copy() | Out-Null
$ws.Range("A4").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$ws.Range("A4").Value = "valid"

# --- B4: "hello@yahoo.com" as a mailto hyperlink (same Hyperlink style as B2/B3) ---
# Create the hyperlink relationship first...
$ws.Hyperlinks.Add($ws.Range("B4"), "mailto:hello@yahoo.com") | Out-Null
# ...then re-apply the exact same formatting used by the other hyperlink
# cells so B4 ends up sharing their cell style, and set the display text.
$ws.Range("B2").Copy() | Out-Null
$ws.Range("B4").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$ws.Range("B4").Value = "hello@yahoo.com"

# --- C4: "qew2" (same plain text style as C2/C3) ---
$ws.Range("C2").Copy() | Out-Null
$ws.Range("C4").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$ws.Range("C4").Value = "qew2"

# --- D4: literal text "TRUE" (NOT a boolean, unlike D2/D3) ---
$ws.Range("A2").Copy() | Out-Null
$ws.Range("D4").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
# Build the literal string "TRUE" via a formula on a scratch cell (typing
# "TRUE" directly into a cell gets auto-coerced to a boolean), then bring
# only the computed value over as text.
$ws.Range("Z1").Formula = "=""TRUE"""
$ws.Range("Z1").Copy() | Out-Null
$ws.Range("D4").PasteSpecial(-4163) | Out-Null   # xlPasteValues
$ws.Range("Z1").Clear() | Out-Null

# --- Update the sheet's active selection ---
$ws.Range("D9").Select() | Out-Null
